$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- Fill in the body placeholder text (shape id=3, "Content Placeholder 2") ---
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$tr.Text = "Build executable`r`r`rOnly execute targets with modified dependencies`rdependency tracking`rsaves lots of time on large projects`rClean all build artifacts"

$tr.Paragraphs(5).IndentLevel = 1
$tr.Paragraphs(6).IndentLevel = 1

# --- Add textbox 1: "$  make" ---
$tb1 = $s.Shapes.AddTextbox(1, (1445341/12700.0), (2458063/12700.0), (2114681/12700.0), (369332/12700.0))
$tb1.Name = "TextBox 4"
$tb1.Fill.ForeColor.RGB = 0
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.TextRange.Text = "`$  make"
$tb1.TextFrame.TextRange.Font.Bold = $true
$tb1.TextFrame.TextRange.Font.Color.RGB = 16777215
$tb1.TextFrame.TextRange.Font.Name = "Courier New"

# --- Add textbox 2: "$  make  clean" ---
$tb2 = $s.Shapes.AddTextbox(1, (1445342/12700.0), (5245506/12700.0), (2114681/12700.0), (369332/12700.0))
$tb2.Name = "TextBox 5"
$tb2.Fill.ForeColor.RGB = 0
$tb2.TextFrame.WordWrap = $false
$tb2.TextFrame.TextRange.Text = "`$  make  clean"
$tb2.TextFrame.TextRange.Font.Bold = $true
$tb2.TextFrame.TextRange.Font.Color.RGB = 16777215
$tb2.TextFrame.TextRange.Font.Name = "Courier New"

# --- Animations ---
$seq = $s.TimeLine.MainSequence

# Effect 1: body placeholder, paragraph 1 (index 0) appears with slide (withEffect)
$e1 = $seq.AddEffect($body, 1, 0, 0)
$e1.Timing.TriggerType = 0
$e1.Exit = $false
$e1.TextRangeStart = 1
$e1.TextRangeLength = 1

$e2 = $seq.AddEffect($tb1, 1, 0, 0)
$e2.Timing.TriggerType = 0

$e3 = $seq.AddEffect($body, 1, 0, 0)
$e4 = $seq.AddEffect($body, 1, 0, 0)
$e5 = $seq.AddEffect($body, 1, 0, 0)

$e6 = $seq.AddEffect($tb2, 1, 0, 0)
